$wb = $excel.ActiveWorkbook

# --- Overview sheet: status text + handoff generate-date bump, E/F/G columns ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-19 12:38:57"
$wsOverview.Range("E:E").ColumnWidth = 16.25
$wsOverview.Range("F:F").ColumnWidth = 16.25

# --- zh-cn sheet: status text + handoff datetime bump, Status column width ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = "Ready for handoff"
$wsZh.Range("H2").Value = "2016-08-19 12:38:53"
$wsZh.Range("C:C").ColumnWidth = 16.25

# --- de-de sheet: status text + handoff datetime bump, Status column width ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = "Ready for handoff"
$wsDe.Range("H2").Value = "2016-08-19 12:38:57"
$wsDe.Range("C:C").ColumnWidth = 16.25
